$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values as text (e.g. "63.021.92", "0.170") so
# that thousand-separator dots and significant trailing zeros are preserved
# exactly as scraped. Force the column back to Text format before writing,
# otherwise Excel's COM layer will auto-coerce plain numeric-looking strings
# (e.g. "147.51") into real numbers and drop formatting like trailing zeros.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '63.097.68'
$ws.Range('E2').Value = '  -0.77%  '
$ws.Range('D3').Value = '2.556.89'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '582.99'
$ws.Range('E5').Value = '  +1.65%  '
$ws.Range('D6').Value = '147.51'
$ws.Range('E6').Value = '  -1.99%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -1.05%  '
$ws.Range('E9').Value = '  -1.07%  '
$ws.Range('E10').Value = '  -4.27%  '
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('D13').Value = '27.25'
$ws.Range('E13').Value = '  -3.50%  '
$ws.Range('D14').Value = '3.012.48'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('D15').Value = '63.003.59'
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('E16').Value = '  -1.06%  '
$ws.Range('D17').Value = '2.543.45'
$ws.Range('E17').Value = '  -1.39%  '
$ws.Range('D18').Value = '11.35'
$ws.Range('E18').Value = '  -3.03%  '
$ws.Range('D19').Value = '336.82'
$ws.Range('E19').Value = '  -1.72%  '
$ws.Range('D20').Value = '4.34'
$ws.Range('E20').Value = '  -1.22%  '
$ws.Range('D21').Value = '6.77'
$ws.Range('E21').Value = '  -2.62%  '
$ws.Range('D23').Value = '65.63'
$ws.Range('E23').Value = '  -0.85%  '
$ws.Range('B24').Value = 'Fetch.AI'
$ws.Range('C24').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D24').Value = '1.63'
$ws.Range('E24').Value = '  +2.25%  '
$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D25').Value = '0.170'
$ws.Range('E25').Value = '  -0.63%  '
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('E28').Value = '  -1.56%  '
$ws.Range('D29').Value = '7.37'
$ws.Range('E29').Value = '  +2.08%  '
$ws.Range('D30').Value = '1.93'
$ws.Range('E30').Value = '  +2.68%  '
$ws.Range('D31').Value = '0.0₃0815'
$ws.Range('E31').Value = '  -3.25%  '
$ws.Range('D32').Value = '177.54'
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('E33').Value = '  -1.96%  '
$ws.Range('D34').Value = '412.47'
$ws.Range('E34').Value = '  -1.41%  '
$ws.Range('D35').Value = '19.16'
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('E36').Value = '  -2.13%  '
$ws.Range('D38').Value = '4.34'
$ws.Range('E38').Value = '  -2.67%  '
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('E40').Value = '  +0.16%  '
$ws.Range('D41').Value = '39.77'
$ws.Range('E41').Value = '  -0.78%  '
$ws.Range('D42').Value = '151.33'
$ws.Range('E42').Value = '  -3.10%  '
$ws.Range('E43').Value = '  -1.53%  '
$ws.Range('D44').Value = '20.92'
$ws.Range('E44').Value = '  -1.82%  '
$ws.Range('D45').Value = '0.0539'
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('E46').Value = '  -1.09%  '
$ws.Range('D47').Value = '0.0969'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('E48').Value = '  +1.65%  '
$ws.Range('D49').Value = '18.39'
$ws.Range('E49').Value = '  -2.52%  '
$ws.Range('E50').Value = '  -7.98%  '
$ws.Range('E51').Value = '  -0.11%  '
